$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    # Force a literal-text value (avoids Excel's automatic date / number
    # coercion for strings such as "2025-05-07" or "000877") and then
    # strip the quote-prefix formatting that Excel would otherwise leave
    # behind, so the cell ends up as a plain shared-string with no style
    # change.
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

# Data rows 2..11 (columns A..H)
# A=Dia(code) B=cliente-date C=quantidade_atipica D=cliente E=id_produto F=produto G=estoque_atualizado H=critico
$data = @(
    @(3, "2025-05-07", 30,  "V V REFEICOES LTDA",                              "000877", "CABO DE ALUMINIO NOBRE 140 CM COM PONTEIRA",             181,  $true),
    @(0, "2025-05-08", 60,  "ASSOCIACAO DOS EMPREENDEDORES DO SMVN",           "000084", "SACO DE LIXO 200L PRETO 0.10 REFORCADO - PCT C/100 UND", -15,  $false),
    @(7, "2025-05-08", 50,  "ASSOCIACAO DOS EMPREENDEDORES DO SMVN",           "000079", "SACO DE LIXO 50L REFORCADO - PCT C/100 UND",              -15,  $false),
    @(9, "2025-05-08", 30,  "CONDOMINIO SOBERANE RESIDENCE, CORPORATE E MALL", "000890", "AROMATIZANTE LIMPADOR PERF CONC COALA ALGODAO 120ML",     25,   $false),
    @(2, "2025-05-09", 10,  "AMAZONIA REFEICOES E SERVICOS LTDA",              "000425", "COADOR DE CAFE EG (EXTRA GRANDE)",                        -3,   $false),
    @(4, "2025-05-09", 15,  "AMAZONIA REFEICOES E SERVICOS LTDA",              "001023", "FILME PVC 30X8X500MT",                                    0,    $true),
    @(5, "2025-05-09", 10,  "AMAZONIA REFEICOES E SERVICOS LTDA",              "000924", "COPO POTE DESCARTAVEL TRANSP 100ML CX/20",                -1,   $false),
    @(8, "2025-05-09", 20,  "A R C DOS SANTOS",                                "000770", "DETERGENTE CLORADO AUDAX GOLD 5L",                        -5,   $true),
    @(1, "2025-05-14", 80,  "MUSASHI DA AMAZONIA LTDA",                        "000307", "DISCO REMOVEDOR PRETO P ENCERADEIRA 510MM",               62,   $false),
    @(6, "2025-05-14", 100, "MUSASHI DA AMAZONIA LTDA",                        "000415", "DETERGENTE DESENGRAX MAX PINE AUDAX 5L",                  163,  $true)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    Set-TextCell $ws.Cells.Item($row, 2) $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    Set-TextCell $ws.Cells.Item($row, 5) $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
    $row++
}

# Row 11 is brand new -- column A carries the bold/bordered "index column"
# style (same as A2:A10). Copy that formatting across from A10 since a
# freshly written cell otherwise comes back with the default style.
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(11, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
